# chore: add dummy data
#
# Adds a "quantity" header in F1 and a first data row (serial number "1",
# book "book1") in row 2, matching the new used range A1:F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell.
$ws.Range("F1").Value = "quantity"

# New dummy data row. The "serial number" column in A1:B1 is stored as text
# (see the numberStoredAsText ignoredError on the sheet), so force A2's
# number format to Text before writing "1" so it round-trips as a string
# instead of being coerced to a numeric value; then restore the cell style
# so no stray formatting is left behind.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "book1"
